$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": append a new data row (38) for date 2020-05-15 (serial 43966),
# pushing the footnote row down from 38 -> 39.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows(38).Insert()
$wsAll.Range("A37:H37").Copy()
$wsAll.Range("A38:H38").PasteSpecial(-4122)
$wsAll.Cells.Item(38, 1).Value = 43966
$wsAll.Cells.Item(38, 2).Value = 281
$wsAll.Cells.Item(38, 3).Value = 277
$wsAll.Cells.Item(38, 4).Value = 64
$wsAll.Cells.Item(38, 5).Value = 54
$wsAll.Cells.Item(38, 6).Value = 10
$wsAll.Cells.Item(38, 7).Value = 11
$wsAll.Cells.Item(38, 8).Value = 202

# ---------------------------------------------------------------------------
# Sheet "kobe": append a new data row (93) for date 2020-05-15 (serial 43966),
# pushing the footnote row down from 93 -> 94.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows(93).Insert()
$wsKobe.Range("A92:J92").Copy()
$wsKobe.Range("A93:J93").PasteSpecial(-4122)
$wsKobe.Cells.Item(93, 1).Value = 43966
$wsKobe.Cells.Item(93, 2).Value = 0
$wsKobe.Cells.Item(93, 3).Value = 2771
$wsKobe.Cells.Item(93, 4).Value = 0
$wsKobe.Cells.Item(93, 5).Value = 281
$wsKobe.Cells.Item(93, 6).Value = 59
$wsKobe.Cells.Item(93, 7).Value = 50
$wsKobe.Cells.Item(93, 8).Value = 9
$wsKobe.Cells.Item(93, 9).Value = 11
$wsKobe.Cells.Item(93, 10).Value = 193

# ---------------------------------------------------------------------------
# Sheet "other": append a new data row (68) for date 2020-05-15 (serial
# 43966), pushing the footnote row down from 68 -> 69.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows(68).Insert()
$wsOther.Range("A67:H67").Copy()
$wsOther.Range("A68:H68").PasteSpecial(-4122)
$wsOther.Cells.Item(68, 1).Value = 43966
$wsOther.Cells.Item(68, 2).Value = 0
$wsOther.Cells.Item(68, 3).Value = 14
$wsOther.Cells.Item(68, 4).Value = 5
$wsOther.Cells.Item(68, 5).Value = 4
$wsOther.Cells.Item(68, 6).Value = 1
$wsOther.Cells.Item(68, 7).Value = 0
$wsOther.Cells.Item(68, 8).Value = 9

# ---------------------------------------------------------------------------
# View state: move the active/selected cell on each of the touched sheets,
# and switch the active workbook tab from "kobe" to "other".
# ---------------------------------------------------------------------------
$wsAll.Range("G42").Select() | Out-Null
$wsKobe.Range("A93").Select() | Out-Null

$wsOther.Activate() | Out-Null
$wsOther.Range("A68").Select() | Out-Null
